$d = $word.ActiveDocument

# 1. Narrative paragraph: "an overall accuracy of 93.3%." -> "an overall accuracy of 95%."
#    The original text " 93.3%" is split across three runs: " 9" | "3.3" | "%".
#    Replace within each run's exact text so run boundaries (and the "%" run) are preserved,
#    mirroring the source edit: " 9" -> " " and "3.3" -> "95".
$d.Content.Find.Execute("overall accuracy of 9", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "overall accuracy of ", 2) | Out-Null
$d.Content.Find.Execute("3.3%.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "95%.", 2) | Out-Null

# 2. Table 1 ("Provided dataset" row), Accuracy (%) cell: "57.2%" -> "85.1"
$d.Content.Find.Execute("57.2%", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "85.1", 2) | Out-Null

# 3. Table 2 ("Custom Dataset (25)" row), Accuracy (%) cell: "93.3" -> "95"
$d.Content.Find.Execute("93.3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "95", 2) | Out-Null
